$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("A4").Value = 75
$wsData.Activate() | Out-Null
$wsData.Range("F1").Select() | Out-Null
